$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: uppercase the language codes in the column header labels
$ws.Range("A1").Value = "instr_msg_EN"
$ws.Range("C1").Value = "instr_msg_ES"
$ws.Range("D1").Value = "instr_msg_FR"

# Swap the image size values for the instructions row (row 2)
$ws.Range("E2").Value = 0.6
$ws.Range("F2").Value = 0.5

# Reset formatting on row 4's size cells back to the default style
$ws.Range("E4:F4").Style = "Normal"

$ws.Range("F2").Select() | Out-Null
